$d = $word.ActiveDocument

# Locate the paragraph that ends with the Nigeria overview text and insert a
# brand-new paragraph right after it containing the Abuja/Lagos sentence.
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphAfter()

$p2 = $d.Paragraphs.Item(2)
$p2.Range.Text = "Abuja is the country" + [char]0x2019 + "s capital located in the North central. Lagos is the main economic hub located in the Western part of the country."
